$d = $word.ActiveDocument
$p15 = $d.Paragraphs.Item(15)
$p24 = $d.Paragraphs.Item(24)
$start = $p15.Range.Start
$end = $p24.Range.End
$r = $d.Range($start, $end)
$r.Delete()

$bm = $d.Bookmarks.Item("_GoBack")
$insPoint = $d.Range($bm.End, $bm.End)
Write-Host ("insPoint Start=" + $insPoint.Start + " End=" + $insPoint.End)
$insPoint.LanguageID = "en-US"
Write-Host ("after set - insPoint Start=" + $insPoint.Start + " End=" + $insPoint.End)
$insPoint.InsertParagraphAfter()
Write-Host "Paragraphs count:" $d.Paragraphs.Count
